$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 962.6949
$ws.Range("I129").Value = 768
$ws.Range("K129").Value = 2304
$ws.Range("M129").Value = 2696
$ws.Range("H137").Value = 1351.42
$ws.Range("I137").Value = 1086.0488
$ws.Range("J137").Value = 2560.3333
$ws.Range("K137").Value = 3258.1464
$ws.Range("L137").Value = 7680.999899999999
$ws.Range("M137").Value = -708.1464000000001
$ws.Range("N137").Value = -12780.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 18997.5
$ws.Range("J55").Value = 18997.5
$ws.Range("L55").Value = 18997.5
$ws.Range("N55").Value = -19627.5
$ws.Range("H61").Value = 3064.3396
$ws.Range("I61").Value = 3227.0222
$ws.Range("J61").Value = 2149.25
$ws.Range("K61").Value = 3227.0222
$ws.Range("L61").Value = 2149.25
$ws.Range("M61").Value = -3015.0222
$ws.Range("N61").Value = -2573.25
$ws.Range("H74").Value = 1285.7941
$ws.Range("I74").Value = 1161.0416
$ws.Range("J74").Value = 1585.2
$ws.Range("K74").Value = 1161.0416
$ws.Range("L74").Value = 1585.2
$ws.Range("M74").Value = -287.0416
$ws.Range("N74").Value = -3333.2
$ws.Range("H77").Value = 1285.7941
$ws.Range("I77").Value = 1161.0416
$ws.Range("J77").Value = 1585.2
$ws.Range("K77").Value = 5805.208000000001
$ws.Range("L77").Value = 7926
$ws.Range("M77").Value = -1437.208000000001
$ws.Range("N77").Value = -16662
$ws.Range("H136").Value = 3064.3396
$ws.Range("I136").Value = 3227.0222
$ws.Range("J136").Value = 2149.25
$ws.Range("K136").Value = 9681.0666
$ws.Range("L136").Value = 6447.75
$ws.Range("M136").Value = -7131.0666
$ws.Range("N136").Value = -11547.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5550.2153
$ws.Range("I31").Value = 1485.9318
$ws.Range("J31").Value = 10659.6
$ws.Range("K31").Value = 1485.9318
$ws.Range("L31").Value = 10659.6
$ws.Range("M31").Value = -1190.9318
$ws.Range("N31").Value = -11249.6
$ws.Range("H34").Value = 5550.2153
$ws.Range("I34").Value = 1485.9318
$ws.Range("J34").Value = 10659.6
$ws.Range("K34").Value = 1485.9318
$ws.Range("L34").Value = 10659.6
$ws.Range("M34").Value = -1283.9318
$ws.Range("N34").Value = -11063.6
$ws.Range("H58").Value = 1331.5106
$ws.Range("I58").Value = 839.3823
$ws.Range("K58").Value = 839.3823
$ws.Range("M58").Value = -636.3823
$ws.Range("H132").Value = 1988.7441
$ws.Range("I132").Value = 1625
$ws.Range("J132").Value = 2602.5625
$ws.Range("K132").Value = 4875
$ws.Range("L132").Value = 7807.6875
$ws.Range("M132").Value = -2345
$ws.Range("N132").Value = -12867.6875
$ws.Range("H136").Value = 1331.5106
$ws.Range("I136").Value = 839.3823
$ws.Range("K136").Value = 2518.1469
$ws.Range("M136").Value = 31.85310000000027
$ws.Range("H141").Value = 34611.152
$ws.Range("J141").Value = 34611.152
$ws.Range("L141").Value = 34611.152
$ws.Range("N141").Value = -44971.152

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 149.83333
$ws.Range("I11").Value = 124.75
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 374.25
$ws.Range("L11").Value = 600
$ws.Range("M11").Value = -234.25
$ws.Range("N11").Value = -880
$ws.Range("H117").Value = 19617920
$ws.Range("I117").Value = 20325.8
$ws.Range("J117").Value = 27783586
$ws.Range("K117").Value = 60977.39999999999
$ws.Range("L117").Value = 83350758
$ws.Range("M117").Value = -57535.39999999999
$ws.Range("N117").Value = -83357642
$ws.Range("H121").Value = 913.9756
$ws.Range("I121").Value = 603.3333
$ws.Range("K121").Value = 1809.9999
$ws.Range("M121").Value = -499.9999
$ws.Range("H122").Value = 3092.7568
$ws.Range("I122").Value = 506.92307
$ws.Range("J122").Value = 4493.4165
$ws.Range("K122").Value = 4562.30763
$ws.Range("L122").Value = 40440.7485
$ws.Range("M122").Value = -2112.30763
$ws.Range("N122").Value = -45340.7485
$ws.Range("H129").Value = 22223372
$ws.Range("I129").Value = 47619624
$ws.Range("J129").Value = 1651.5
$ws.Range("K129").Value = 142858872
$ws.Range("L129").Value = 4954.5
$ws.Range("M129").Value = -142853872
$ws.Range("N129").Value = -14954.5
$ws.Range("H138").Value = 29421594
$ws.Range("I138").Value = 35725576
$ws.Range("J138").Value = 3021
$ws.Range("K138").Value = 107176728
$ws.Range("L138").Value = 9063
$ws.Range("M138").Value = -107171588
$ws.Range("N138").Value = -19343
$ws.Range("H141").Value = 9767.4
$ws.Range("I141").Value = 9875.916999999999
$ws.Range("J141").Value = 9333.333000000001
$ws.Range("K141").Value = 29627.751
$ws.Range("L141").Value = 27999.999
$ws.Range("M141").Value = -24447.751
$ws.Range("N141").Value = -38359.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1698.6923
$ws.Range("I132").Value = 1392.0513
$ws.Range("J132").Value = 2618.6155
$ws.Range("K132").Value = 4176.1539
$ws.Range("L132").Value = 7855.8465
$ws.Range("M132").Value = -1646.1539
$ws.Range("N132").Value = -12915.8465
$ws.Range("H134").Value = 14707.143
$ws.Range("J134").Value = 14707.143
$ws.Range("L134").Value = 44121.429
$ws.Range("N134").Value = -49191.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1808.0714
$ws.Range("I7").Value = 1627.3
$ws.Range("J7").Value = 2260
$ws.Range("K7").Value = 1627.3
$ws.Range("L7").Value = 2260
$ws.Range("M7").Value = -1515.3
$ws.Range("N7").Value = -2484
$ws.Range("H126").Value = 1808.0714
$ws.Range("I126").Value = 1627.3
$ws.Range("J126").Value = 2260
$ws.Range("K126").Value = 4881.9
$ws.Range("L126").Value = 6780
$ws.Range("M126").Value = -2411.9
$ws.Range("N126").Value = -11720
$ws.Range("H132").Value = 16703597
$ws.Range("I132").Value = 23238498
$ws.Range("J132").Value = 3295.4443
$ws.Range("K132").Value = 69715494
$ws.Range("L132").Value = 9886.332900000001
$ws.Range("M132").Value = -69712964
$ws.Range("N132").Value = -14946.3329
$ws.Range("H136").Value = 5075.783
$ws.Range("I136").Value = 3502.76
$ws.Range("J136").Value = 12940.9
$ws.Range("K136").Value = 10508.28
$ws.Range("L136").Value = 38822.7
$ws.Range("M136").Value = -7958.280000000001
$ws.Range("N136").Value = -43922.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 602.5333000000001
$ws.Range("I126").Value = 467
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 1401
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = 1069
$ws.Range("N126").Value = -12440
$ws.Range("H132").Value = 20899.56
$ws.Range("I132").Value = 27036.078
$ws.Range("J132").Value = 1467.25
$ws.Range("K132").Value = 81108.234
$ws.Range("L132").Value = 4401.75
$ws.Range("M132").Value = -78578.234
$ws.Range("N132").Value = -9461.75
$ws.Range("H136").Value = 5815751
$ws.Range("I136").Value = 1828.8197
$ws.Range("J136").Value = 20001720
$ws.Range("K136").Value = 5486.4591
$ws.Range("L136").Value = 20001720
$ws.Range("M136").Value = -2936.4591
$ws.Range("N136").Value = -60010260
